# Update gh-pages output data: refresh "want to go" counts (column F) and
# mark the cancelled "苏州·明日方舟ONLY#2024~佑桑柔" event (column C) on
# both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetExhibit = $wb.Worksheets.Item("展览")
$sheetAll     = $wb.Worksheets.Item("全部类型")

# --- 展览 (rows as in the source sheet) ---
$sheetExhibit.Range("F3").Value  = 12944
$sheetExhibit.Range("F5").Value  = 78
$sheetExhibit.Range("F6").Value  = 90
$sheetExhibit.Range("F10").Value = 12907
$sheetExhibit.Range("F12").Value = 41
$sheetExhibit.Range("F13").Value = 8699
$sheetExhibit.Range("F14").Value = 7703
$sheetExhibit.Range("F16").Value = 113
$sheetExhibit.Range("C17").Value = "苏州·明日方舟ONLY#2024~佑桑柔（取消）"
$sheetExhibit.Range("F19").Value = 986
$sheetExhibit.Range("F24").Value = 20

# --- 全部类型 (same rows, shifted by one because of the extra 演出 entry) ---
$sheetAll.Range("F4").Value  = 12944
$sheetAll.Range("F6").Value  = 78
$sheetAll.Range("F7").Value  = 90
$sheetAll.Range("F11").Value = 12907
$sheetAll.Range("F13").Value = 41
$sheetAll.Range("F14").Value = 8699
$sheetAll.Range("F15").Value = 7703
$sheetAll.Range("F17").Value = 113
$sheetAll.Range("C18").Value = "苏州·明日方舟ONLY#2024~佑桑柔（取消）"
$sheetAll.Range("F20").Value = 986
$sheetAll.Range("F26").Value = 20
